# feat: Tiled map class implementation
#
# Adds a new "Tile Types" worksheet (placed after the existing Animals /
# DayNightCycle sheets) describing, per tile id, the boolean flags consumed
# by the tiled-map system: IsStructure, IsWalkable, IsAnimalHabitat,
# CanSpawnHumans, CanDispawnHumans.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet so it lands at the
# end of the tab strip: Animals, DayNightCycle, Tile Types.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tile Types"

# Header row
$ws.Range("A1").Value = "IDS"
$ws.Range("B1").Value = "Is Structure"
$ws.Range("C1").Value = "Is Walkable"
$ws.Range("D1").Value = "Is Animal Habitat"
$ws.Range("E1").Value = "Can Spawn Humans"
$ws.Range("F1").Value = "Can Dispawn Humans"

# Row 2 - Road
$ws.Range("A2").Value = "Road"
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = $false

# Row 3 - Jail Walls
$ws.Range("A3").Value = "Jail Walls"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = $false

# Row 4 - Jail Habitat
$ws.Range("A4").Value = "Jail Habitat"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = $false

# Row 5 - Structure
$ws.Range("A5").Value = "Structure"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = $false

# Row 6 - Humans Entry
$ws.Range("A6").Value = "Humans Entry"
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = $false

# Row 7 - Humans Exit
$ws.Range("A7").Value = "Humans Exit"
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $true

# Row 8 - Grass
$ws.Range("A8").Value = "Grass"
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $false

# Match the sheet-wide default font/formatting (Arial, theme text color,
# minor-scheme) that a freshly created sheet carries.
$ws.Range("A1:F8").Font.ThemeColor = 1
